$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (rows 2-51); rows 45/46 also swap Coin/Link (TheGraph <-> USDe)

$ws.Cells.Item(2, 4).Value = '62.805.14'
$ws.Cells.Item(2, 5).Value = '  +2.81%  '

$ws.Cells.Item(3, 4).Value = '3.039.77'
$ws.Cells.Item(3, 5).Value = '  +1.40%  '

$ws.Cells.Item(4, 5).Value = '  +0.29%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '548.32'
$ws.Cells.Item(5, 5).Value = '  +3.56%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '137.05'
$ws.Cells.Item(6, 5).Value = '  +3.72%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '1.00'
$ws.Cells.Item(7, 5).Value = '  +0.16%  '

$ws.Cells.Item(8, 4).Value = '3.035.33'
$ws.Cells.Item(8, 5).Value = '  +1.32%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.497'
$ws.Cells.Item(9, 5).Value = '  +1.25%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '6.21'
$ws.Cells.Item(10, 5).Value = '  +0.91%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.148'
$ws.Cells.Item(11, 5).Value = '  -0.82%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.446'
$ws.Cells.Item(12, 5).Value = '  +0.76%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.0000224'
$ws.Cells.Item(13, 5).Value = '  +2.51%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '34.40'
$ws.Cells.Item(14, 5).Value = '  +2.08%  '

$ws.Cells.Item(15, 4).Value = '3.544.78'
$ws.Cells.Item(15, 5).Value = '  +2.02%  '

$ws.Cells.Item(16, 4).Value = '62.936.34'
$ws.Cells.Item(16, 5).Value = '  +3.14%  '

$ws.Cells.Item(17, 4).Value = '3.049.56'
$ws.Cells.Item(17, 5).Value = '  +2.05%  '

$ws.Cells.Item(18, 5).Value = '  -1.67%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '6.66'
$ws.Cells.Item(19, 5).Value = '  +1.75%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '477.18'
$ws.Cells.Item(20, 5).Value = '  +2.69%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '13.48'
$ws.Cells.Item(21, 5).Value = '  +2.24%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '0.668'
$ws.Cells.Item(22, 5).Value = '  -0.60%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '7.11'
$ws.Cells.Item(23, 5).Value = '  +3.06%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '80.50'
$ws.Cells.Item(24, 5).Value = '  +1.98%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '12.33'
$ws.Cells.Item(25, 5).Value = '  +3.81%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.998'
$ws.Cells.Item(26, 5).Value = '  -0.23%  '

$ws.Cells.Item(27, 5).Value = '  +3.00%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '7.78'
$ws.Cells.Item(28, 5).Value = '  +0.35%  '

$ws.Cells.Item(29, 5).Value = '  +0.30%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '1.96'
$ws.Cells.Item(30, 5).Value = '  +4.88%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '25.73'
$ws.Cells.Item(31, 5).Value = '  +1.32%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '1.14'
$ws.Cells.Item(32, 5).Value = '  -0.25%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '2.38'
$ws.Cells.Item(33, 5).Value = '  +5.58%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '5.61'
$ws.Cells.Item(34, 5).Value = '  +3.65%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '54.99'
$ws.Cells.Item(35, 5).Value = '  -0.44%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '5.90'
$ws.Cells.Item(36, 5).Value = '  +1.38%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '458.06'
$ws.Cells.Item(37, 5).Value = '  -0.20%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.0805'
$ws.Cells.Item(38, 5).Value = '  +2.92%  '

$ws.Cells.Item(39, 4).Value = '3.099.35'
$ws.Cells.Item(39, 5).Value = '  -2.64%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.0390'
$ws.Cells.Item(40, 5).Value = '  +2.00%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.116'
$ws.Cells.Item(41, 5).Value = '  -0.80%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '8.16'
$ws.Cells.Item(42, 5).Value = '  +0.98%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '2.51'
$ws.Cells.Item(43, 5).Value = '  +2.70%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '27.62'
$ws.Cells.Item(44, 5).Value = '  +5.17%  '

$ws.Cells.Item(45, 2).Value = 'USDe'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '1.00'
$ws.Cells.Item(45, 5).Value = '  -0.12%  '

$ws.Cells.Item(46, 2).Value = 'TheGraph'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.249'
$ws.Cells.Item(46, 5).Value = '  +1.55%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '2.01'
$ws.Cells.Item(47, 5).Value = '  +2.03%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.108'
$ws.Cells.Item(48, 5).Value = '  +0.56%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '115.98'
$ws.Cells.Item(49, 5).Value = '  -2.76%  '

$ws.Cells.Item(50, 4).Value = '0.0₃0499'
$ws.Cells.Item(50, 5).Value = '  +0.91%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '2.04'
$ws.Cells.Item(51, 5).Value = '  +3.25%  '
